$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = 'Configuration setup'
$ws.Range("B23").Value = 'Read all the configurations'
$ws.Range("C23").Value = 'Pass'
$ws.Range("D23").Value = '21-01-2021 05:11:30 PM'

$ws.Range("A24").Value = 'Login to Portal '
$ws.Range("B24").Value = 'Logged in sucessfully '
$ws.Range("C24").Value = 'Pass'
$ws.Range("D24").Value = '21-01-2021 05:11:56 PM'

$ws.Range("A25").Value = 'Navigation to INT'
$ws.Range("B25").Value = 'Successfully Navigated to International '
$ws.Range("C25").Value = 'Pass'
$ws.Range("D25").Value = '21-01-2021 05:12:11 PM'

$ws.Range("A26").Value = 'Search of created WO '
$ws.Range("B26").Value = 'Search is Working'
$ws.Range("C26").Value = 'Pass'
$ws.Range("D26").Value = '21-01-2021 05:12:23 PM'

$ws.Range("A27").Value = 'House format (HF) medialevel'
$ws.Range("B27").Value = 'House format (HF) medialevel is completed'
$ws.Range("C27").Value = 'Pass'
$ws.Range("D27").Value = '21-01-2021 05:12:38 PM'

$ws.Range("A28").Value = 'Configuration setup'
$ws.Range("B28").Value = 'Read all the configurations'
$ws.Range("C28").Value = 'Pass'
$ws.Range("D28").Value = '21-01-2021 05:16:11 PM'

$ws.Range("A29").Value = 'Login to Portal '
$ws.Range("B29").Value = 'Logged in sucessfully '
$ws.Range("C29").Value = 'Pass'
$ws.Range("D29").Value = '21-01-2021 05:16:40 PM'

$ws.Range("A30").Value = 'Navigation to INT'
$ws.Range("B30").Value = 'Successfully Navigated to International '
$ws.Range("C30").Value = 'Pass'
$ws.Range("D30").Value = '21-01-2021 05:16:54 PM'

$ws.Range("A31").Value = 'Search of created WO '
$ws.Range("B31").Value = 'Search is Working'
$ws.Range("C31").Value = 'Pass'
$ws.Range("D31").Value = '21-01-2021 05:17:07 PM'

$ws.Range("A32").Value = 'House format (HF) medialevel'
$ws.Range("B32").Value = 'House format (HF) medialevel is completed'
$ws.Range("C32").Value = 'Pass'
$ws.Range("D32").Value = '21-01-2021 05:17:21 PM'

$ws.Range("A33").Value = 'Configuration setup'
$ws.Range("B33").Value = 'Read all the configurations'
$ws.Range("C33").Value = 'Pass'
$ws.Range("D33").Value = '21-01-2021 05:18:15 PM'

$ws.Range("A34").Value = 'Login to Portal '
$ws.Range("B34").Value = 'Logged in sucessfully '
$ws.Range("C34").Value = 'Pass'
$ws.Range("D34").Value = '21-01-2021 05:18:43 PM'

$ws.Range("A35").Value = 'Navigation to INT'
$ws.Range("B35").Value = 'Successfully Navigated to International '
$ws.Range("C35").Value = 'Pass'
$ws.Range("D35").Value = '21-01-2021 05:18:58 PM'

$ws.Range("A36").Value = 'Search of created WO '
$ws.Range("B36").Value = 'Search is Working'
$ws.Range("C36").Value = 'Pass'
$ws.Range("D36").Value = '21-01-2021 05:19:10 PM'

$ws.Range("A37").Value = 'House format (HF) medialevel'
$ws.Range("B37").Value = 'House format (HF) medialevel is completed'
$ws.Range("C37").Value = 'Pass'
$ws.Range("D37").Value = '21-01-2021 05:19:25 PM'

$ws.Range("A38").Value = 'AQC overRide is clicked '
$ws.Range("B38").Value = 'AQC over Ride performed '
$ws.Range("C38").Value = 'Pass'
$ws.Range("D38").Value = '21-01-2021 05:19:40 PM'

$ws.Range("A39").Value = 'AQC overRide pass is enabled '
$ws.Range("B39").Value = 'AQC over ride Pass'
$ws.Range("C39").Value = 'Fail'
$ws.Range("D39").Value = '21-01-2021 05:19:45 PM'

$ws.Range("A40").Value = 'House format (HF) medialevel'
$ws.Range("B40").Value = 'House format (HF) medialevel is completed'
$ws.Range("C40").Value = 'Pass'
$ws.Range("D40").Value = '21-01-2021 05:20:16 PM'

$ws.Range("A41").Value = 'Configuration setup'
$ws.Range("B41").Value = 'Read all the configurations'
$ws.Range("C41").Value = 'Pass'
$ws.Range("D41").Value = '21-01-2021 05:28:26 PM'

$ws.Range("A42").Value = 'Login to Portal '
$ws.Range("B42").Value = 'Logged in sucessfully '
$ws.Range("C42").Value = 'Pass'
$ws.Range("D42").Value = '21-01-2021 05:28:54 PM'

$ws.Range("A43").Value = 'Navigation to INT'
$ws.Range("B43").Value = 'Successfully Navigated to International '
$ws.Range("C43").Value = 'Pass'
$ws.Range("D43").Value = '21-01-2021 05:29:09 PM'

$ws.Range("A44").Value = 'Search of created WO '
$ws.Range("B44").Value = 'Search is Working'
$ws.Range("C44").Value = 'Pass'
$ws.Range("D44").Value = '21-01-2021 05:29:21 PM'

$ws.Range("A45").Value = 'House format (HF) medialevel'
$ws.Range("B45").Value = 'House format (HF) medialevel is completed'
$ws.Range("C45").Value = 'Pass'
$ws.Range("D45").Value = '21-01-2021 05:29:36 PM'

$ws.Range("A46").Value = 'Configuration setup'
$ws.Range("B46").Value = 'Read all the configurations'
$ws.Range("C46").Value = 'Pass'
$ws.Range("D46").Value = '22-01-2021 12:39:04 PM'

$ws.Range("A47").Value = 'Login to Portal '
$ws.Range("B47").Value = 'Logged in sucessfully '
$ws.Range("C47").Value = 'Pass'
$ws.Range("D47").Value = '22-01-2021 12:39:44 PM'

$ws.Range("A48").Value = 'Navigation to INT'
$ws.Range("B48").Value = 'Successfully Navigated to International '
$ws.Range("C48").Value = 'Pass'
$ws.Range("D48").Value = '22-01-2021 12:39:59 PM'

$ws.Range("A49").Value = 'Search of created WO '
$ws.Range("B49").Value = 'Search is Working'
$ws.Range("C49").Value = 'Pass'
$ws.Range("D49").Value = '22-01-2021 12:40:11 PM'

$ws.Range("A50").Value = 'Configuration setup'
$ws.Range("B50").Value = 'Read all the configurations'
$ws.Range("C50").Value = 'Pass'
$ws.Range("D50").Value = '22-01-2021 12:59:50 PM'

$ws.Range("A51").Value = 'Login to Portal '
$ws.Range("B51").Value = 'Logged in sucessfully '
$ws.Range("C51").Value = 'Pass'
$ws.Range("D51").Value = '22-01-2021 01:00:18 PM'

$ws.Range("A52").Value = 'Configuration setup'
$ws.Range("B52").Value = 'Read all the configurations'
$ws.Range("C52").Value = 'Pass'
$ws.Range("D52").Value = '22-01-2021 03:55:02 PM'

$ws.Range("A53").Value = 'Login to Portal '
$ws.Range("B53").Value = 'Logged in sucessfully '
$ws.Range("C53").Value = 'Pass'
$ws.Range("D53").Value = '22-01-2021 03:55:30 PM'

$ws.Range("A54").Value = 'Navigation to INT'
$ws.Range("B54").Value = 'Successfully Navigated to International '
$ws.Range("C54").Value = 'Pass'
$ws.Range("D54").Value = '22-01-2021 03:55:58 PM'

$ws.Range("A55").Value = 'Search of created WO '
$ws.Range("B55").Value = 'Search is Working'
$ws.Range("C55").Value = 'Pass'
$ws.Range("D55").Value = '22-01-2021 03:56:11 PM'

$ws.Range("A56").Value = 'Ingest status in DM level'
$ws.Range("B56").Value = 'Ingest status in DM level is completed'
$ws.Range("C56").Value = 'Pass'
$ws.Range("D56").Value = '22-01-2021 03:56:25 PM'

$ws.Range("A57").Value = 'AQC overRide is clicked '
$ws.Range("B57").Value = 'AQC over Ride performed '
$ws.Range("C57").Value = 'Pass'
$ws.Range("D57").Value = '22-01-2021 03:56:40 PM'

$ws.Range("A58").Value = 'AQC overRide pass is enabled '
$ws.Range("B58").Value = 'AQC over ride Pass'
$ws.Range("C58").Value = 'Fail'
$ws.Range("D58").Value = '22-01-2021 03:56:45 PM'

$ws.Range("A59").Value = 'Configuration setup'
$ws.Range("B59").Value = 'Read all the configurations'
$ws.Range("C59").Value = 'Pass'
$ws.Range("D59").Value = '22-01-2021 04:01:25 PM'

$ws.Range("A60").Value = 'Login to Portal '
$ws.Range("B60").Value = 'Logged in sucessfully '
$ws.Range("C60").Value = 'Pass'
$ws.Range("D60").Value = '22-01-2021 04:01:53 PM'

$ws.Range("A61").Value = 'Navigation to INT'
$ws.Range("B61").Value = 'Successfully Navigated to International '
$ws.Range("C61").Value = 'Pass'
$ws.Range("D61").Value = '22-01-2021 04:02:07 PM'

$ws.Range("A62").Value = 'Search of created WO '
$ws.Range("B62").Value = 'Search is Working'
$ws.Range("C62").Value = 'Pass'
$ws.Range("D62").Value = '22-01-2021 04:02:20 PM'

$ws.Range("A63").Value = 'Ingest status in DM level'
$ws.Range("B63").Value = 'Ingest status in DM level is completed'
$ws.Range("C63").Value = 'Pass'
$ws.Range("D63").Value = '22-01-2021 04:02:34 PM'

$ws.Range("A64").Value = 'AQC overRide is clicked '
$ws.Range("B64").Value = 'AQC over Ride performed '
$ws.Range("C64").Value = 'Pass'
$ws.Range("D64").Value = '22-01-2021 04:02:49 PM'

$ws.Range("A65").Value = 'AQC overRide pass is enabled '
$ws.Range("B65").Value = 'AQC over ride Pass'
$ws.Range("C65").Value = 'Fail'
$ws.Range("D65").Value = '22-01-2021 04:02:54 PM'

$ws.Range("D66").Value = 'DateTime.Now.Date.ToString("dd-MM-yyyy")'

$ws.Range("A67").Value = 'Configuration setup'
$ws.Range("B67").Value = 'Read all the configurations'
$ws.Range("C67").Value = 'Pass'
$ws.Range("D67").Value = '22-01-2021 04:14:54 PM'

$ws.Range("A68").Value = 'Login to Portal '
$ws.Range("B68").Value = 'Logged in sucessfully '
$ws.Range("C68").Value = 'Pass'
$ws.Range("D68").Value = '22-01-2021 04:15:24 PM'
